$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so values are preserved as literal text
$textCells = @('D5', 'D8', 'D9', 'D11', 'D13', 'D14', 'D19', 'D20', 'D22', 'D23', 'D24', 'D26', 'D27', 'D29', 'D30', 'D32', 'D33', 'D34', 'D38', 'D39', 'D40', 'D42', 'D43', 'D44', 'D45', 'D49', 'D50')
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2: Bitcoin
$ws.Range('D2').Value = '36.970.00'
$ws.Range('E2').Value = '  -1.68%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.017.44'
$ws.Range('E3').Value = '  -3.25%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.05%  '

# Row 5: BNB
$ws.Range('D5').Value = '226.56'
$ws.Range('E5').Value = '  -3.07%  '

# Row 6: XRP
$ws.Range('E6').Value = '  -4.33%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.06%  '

# Row 8: Solana
$ws.Range('D8').Value = '54.95'
$ws.Range('E8').Value = '  -5.44%  '

# Row 9: Cardano
$ws.Range('D9').Value = '0.380'
$ws.Range('E9').Value = '  -2.98%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  +1.51%  '

# Row 11: TRON
$ws.Range('D11').Value = '0.104'
$ws.Range('E11').Value = '  -3.95%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '2.316.58'
$ws.Range('E12').Value = '  -3.27%  '

# Row 13: Chainlink
$ws.Range('D13').Value = '14.26'
$ws.Range('E13').Value = '  -5.73%  '

# Row 14: Avalanche
$ws.Range('D14').Value = '20.59'
$ws.Range('E14').Value = '  -2.62%  '

# Row 15: Polygon
$ws.Range('E15').Value = '  -4.02%  '

# Row 16: Polkadot
$ws.Range('E16').Value = '  -3.97%  '

# Row 17: WrappedEther
$ws.Range('D17').Value = '2.027.39'
$ws.Range('E17').Value = '  -2.75%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '36.871.55'
$ws.Range('E18').Value = '  -1.89%  '

# Row 19: Uniswap
$ws.Range('D19').Value = '6.04'
$ws.Range('E19').Value = '  -0.23%  '

# Row 20: Litecoin
$ws.Range('D20').Value = '68.81'
$ws.Range('E20').Value = '  -2.86%  '

# Row 21: ShibaInu
$ws.Range('D21').Value = '0.0₃0827'
$ws.Range('E21').Value = '  -0.79%  '

# Row 22: BitcoinCash
$ws.Range('D22').Value = '226.11'
$ws.Range('E22').Value = '  -1.39%  '

# Row 23: Dai
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.12%  '

# Row 24: Toncoin
$ws.Range('D24').Value = '2.41'
$ws.Range('E24').Value = '  +1.83%  '

# Row 25: PancakeSwap
$ws.Range('E25').Value = '  -4.72%  '

# Row 26: Monero
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '167.00'
$ws.Range('E26').Value = '  -2.17%  '

# Row 27: Cosmos
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = '9.30'
$ws.Range('E27').Value = '  -4.25%  '

# Row 28: Kaspa
$ws.Range('E28').Value = '  -4.62%  '

# Row 29: EthereumClassic
$ws.Range('D29').Value = '18.76'
$ws.Range('E29').Value = '  -4.13%  '

# Row 30: ImmutableX
$ws.Range('D30').Value = '1.33'
$ws.Range('E30').Value = '  -3.59%  '

# Row 31: Stellar
$ws.Range('E31').Value = '  -4.45%  '

# Row 32: Filecoin
$ws.Range('D32').Value = '4.47'
$ws.Range('E32').Value = '  -4.07%  '

# Row 33: Hedera
$ws.Range('D33').Value = '0.0611'
$ws.Range('E33').Value = '  -4.24%  '

# Row 34: InternetComputer(DFINITY)
$ws.Range('D34').Value = '4.43'
$ws.Range('E34').Value = '  -5.08%  '

# Row 35: LidoDAOToken
$ws.Range('E35').Value = '  -4.91%  '

# Row 36: WEMIXToken
$ws.Range('E36').Value = '  +0.64%  '

# Row 37: BinanceUSD
$ws.Range('E37').Value = '  +0.21%  '

# Row 38: RenderToken
$ws.Range('D38').Value = '3.17'
$ws.Range('E38').Value = '  -4.84%  '

# Row 39: THORChain
$ws.Range('D39').Value = '5.42'
$ws.Range('E39').Value = '  +0.67%  '

# Row 40: VeChain
$ws.Range('B40').Value = 'VeChain'
$ws.Range('C40').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D40').Value = '0.0219'
$ws.Range('E40').Value = '  -5.79%  '

# Row 41: Maker
$ws.Range('B41').Value = 'Maker'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D41').Value = '1.495.19'
$ws.Range('E41').Value = '  +2.40%  '

# Row 42: InjectiveProtocol
$ws.Range('D42').Value = '16.91'
$ws.Range('E42').Value = '  +0.21%  '

# Row 43: Cronos
$ws.Range('B43').Value = 'Cronos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D43').Value = '0.0930'
$ws.Range('E43').Value = '  -3.16%  '

# Row 44: Aave
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').Value = '94.93'
$ws.Range('E44').Value = '  -5.78%  '

# Row 45: HuobiToken
$ws.Range('D45').Value = '2.82'
$ws.Range('E45').Value = '  -2.61%  '

# Row 46: TrustWalletToken
$ws.Range('E46').Value = '  -5.43%  '

# Row 47: FraxShare
$ws.Range('E47').Value = '  -0.15%  '

# Row 49: MXToken
$ws.Range('D49').Value = '2.91'
$ws.Range('E49').Value = '  -1.68%  '

# Row 50: FTXToken
$ws.Range('D50').Value = '3.67'
$ws.Range('E50').Value = '  -8.08%  '

# Row 51: RocketPoolETH
$ws.Range('D51').Value = '2.206.55'
$ws.Range('E51').Value = '  -3.18%  '
